# Added periods to exclude from GAGE
# Adds six new site sheets (MHD, THD, RPB, SMO, ZEP, CMN) after the existing
# "CGO" sheet. Each new sheet carries over just the header block (rows 1-8)
# from CGO - the commented species/notes rows plus the bold column-heading
# row - without the CGO-specific exclusion-period data rows (9-11) or the
# cell comments attached to them.

$wb = $excel.ActiveWorkbook
$cgo = $wb.Worksheets.Item("CGO")

$newNames = @("MHD", "THD", "RPB", "SMO", "ZEP", "CMN")

foreach ($name in $newNames) {
    # Duplicate CGO so the new sheet starts with identical formatting
    # (column widths, bold header row, number formats, etc.)
    $cgo.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
    $newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet.Name = $name

    # Remove the two cell comments that lived on the data rows we are
    # about to delete.
    $newSheet.Range("E9").Comment.Delete()
    $newSheet.Range("C11").Comment.Delete()

    # Drop the CGO-specific instrument-period rows; only the shared header
    # rows (1-8) should remain on the new site sheets.
    $newSheet.Range("A9:M11").Delete()

    # Match the default, unselected view state of a freshly duplicated sheet.
    $newSheet.Range("D13").Select()
}

# Leave the original CGO sheet as the active tab/selection, as before.
$cgo.Activate()
$cgo.Range("H17").Select()
